# Reisezeitmatrix Haltestelle zu Haltestelle
# Adds a machine-readable header row ("from_stop" / "to_stop" / "minutes")
# above the existing (translated) header row, fills in sample travel-time
# data, and adds a data validation rule for the new numeric id columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row 1 above the current header ------------------------
# (the old header "Von Haltestelle (Nr)" / "Zu Haltestelle (Nr)" /
#  "Reisezeit in Minuten" moves down to row 2, together with the
#  dataValidation ranges that referenced it)
$ws.Rows("1:1").Insert()

# --- 2. Technical column headers in the new row 1 ---------------------------
$ws.Range("A1").Value = "from_stop"
$ws.Range("B1").Value = "to_stop"
$ws.Range("C1").Value = "minutes"

# Match the formatting of the (former/ now row 2) header row
$ws.Range("A2:C2").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

# --- 3. Sample data rows 3-10 -----------------------------------------------
$data = @(
    @(3, 3, 10),
    @(5, 3, 20),
    @(77, 3, 30),
    @(99, 3, 40),
    @(3, 5, 50),
    @(5, 5, 60),
    @(77, 5, 70),
    @(99, 5, 80)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- 4. Data validation ------------------------------------------------------
# The existing validations (decimal on column C, whole number on A/B) were
# shifted down to row 2+ by the row insert above. Shrink them back to just
# the header row (row 2) and (re)create a dedicated "whole number between
# 0 and 9999999999" validation for the actual id data starting at row 3.
$ws.Range("A3:B1048576").Validation.Delete()
$ws.Range("A3:B1048576").Validation.Add(1, 1, 1, 0, 9999999999)

# --- 5. Selection -------------------------------------------------------------
$ws.Range("C1").Select()
